$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This commit regenerates the handback-status report with a fresh run's
# file GUIDs and timestamps:
#   383509dd-8e05-4b7e-8273-1d9fa6c3dfe6  ->  21d41568-f00e-4aaf-90d5-3c89a0e5ceb8
#   42647d57-8228-4722-a6e3-4fd76a0d03a6  ->  fffff0741068-fb15-4787-8417-99839806d122
# plus new xlf content-hashes and new handoff/handback timestamps.
# ---------------------------------------------------------------------------

$oldGuid1 = "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6"
$newGuid1 = "21d41568-f00e-4aaf-90d5-3c89a0e5ceb8"
$oldGuid2 = "42647d57-8228-4722-a6e3-4fd76a0d03a6"
$newGuid2 = "fffff0741068-fb15-4787-8417-99839806d122"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"
$newMdPath1 = "e2e\$newGuid1.md"
$newMdPath2 = "e2e\$newGuid2.md"

$newZhXlf = "$newGuid1.885cd247df98f70f27b41408080dead45342786e.zh-cn.xlf"
$newDeXlf = "$newGuid1.885cd247df98f70f27b41408080dead45342786e.de-de.xlf"

$newOverviewDate = "2016-08-24 17:05:31"
$newZhHandoffDate = "2016-08-24 17:05:26"
$newZhHandbackDate = "2016-08-24 17:05:43"
$newDeHandoffDate = "2016-08-24 17:05:31"
$newDeHandbackDate = "2016-08-24 17:05:51"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2,1).Value = $newMd1
$wsOverview.Cells.Item(2,2).Value = $newMdPath1
$wsOverview.Cells.Item(2,7).Value = $newOverviewDate

$wsOverview.Cells.Item(3,1).Value = $newMd2
$wsOverview.Cells.Item(3,2).Value = $newMdPath2
$wsOverview.Cells.Item(3,7).Value = $newOverviewDate

# refresh hyperlinks on B2/B3 so their display text + target follow the new names
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B3").Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/$newMd1", "", "", $newMdPath1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/$newMd2", "", "", $newMdPath2)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(2,1).Value = $newMd1
$wsZh.Cells.Item(2,7).Value = $newZhXlf
$wsZh.Cells.Item(2,8).Value = $newZhHandoffDate
$wsZh.Cells.Item(2,9).Value = $newMd1
$wsZh.Cells.Item(2,10).Value = $newZhXlf
$wsZh.Cells.Item(2,11).Value = $newZhHandbackDate

$wsZh.Cells.Item(3,1).Value = $newMd2
$wsZh.Cells.Item(3,7).Value = $newZhXlf
$wsZh.Cells.Item(3,8).Value = $newZhHandoffDate
$wsZh.Cells.Item(3,9).Value = $newMd2
$wsZh.Cells.Item(3,10).Value = $newZhXlf
$wsZh.Cells.Item(3,11).Value = $newZhHandbackDate

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("A3").Hyperlinks.Delete()
$wsZh.Range("I3").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/$newMd1", "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0e85904641765d8d04647f364d5b04a0c5bee83e/e2e/$newMd1", "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/$newMd2", "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0e85904641765d8d04647f364d5b04a0c5bee83e/e2e/$newMd2", "", "", $newMd2)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(2,1).Value = $newMd1
$wsDe.Cells.Item(2,7).Value = $newDeXlf
$wsDe.Cells.Item(2,8).Value = $newDeHandoffDate
$wsDe.Cells.Item(2,9).Value = $newMd1
$wsDe.Cells.Item(2,10).Value = $newDeXlf
$wsDe.Cells.Item(2,11).Value = $newDeHandbackDate

$wsDe.Cells.Item(3,1).Value = $newMd2
$wsDe.Cells.Item(3,7).Value = $newDeXlf
$wsDe.Cells.Item(3,8).Value = $newDeHandoffDate
$wsDe.Cells.Item(3,9).Value = $newMd2
$wsDe.Cells.Item(3,10).Value = $newDeXlf
$wsDe.Cells.Item(3,11).Value = $newDeHandbackDate

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("A3").Hyperlinks.Delete()
$wsDe.Range("I3").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/$newMd1", "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f79ddf06c7e2584bae9afd6c2a47240203fd580e/e2e/$newMd1", "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/$newMd2", "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f79ddf06c7e2584bae9afd6c2a47240203fd580e/e2e/$newMd2", "", "", $newMd2)
